$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 566, pushing existing rows 566-612 down to 568-614
$ws.Range("A566:R567").EntireRow.Insert()

# Populate the first new row (566) - Pimiento, Zafiro rojo
$ws.Cells.Item(566, 1).Value = 5
$ws.Cells.Item(566, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(566, 3).Value = "Maule"
$ws.Cells.Item(566, 4).Value = 44783
$ws.Cells.Item(566, 5).Value = 7
$ws.Cells.Item(566, 6).Value = 100112002
$ws.Cells.Item(566, 7).Value = "Pimiento"
$ws.Cells.Item(566, 8).Value = "Zafiro rojo"
$ws.Cells.Item(566, 9).Value = "Primera"
$ws.Cells.Item(566, 10).Value = 300
$ws.Cells.Item(566, 11).Value = 25000
$ws.Cells.Item(566, 12).Value = 25000
$ws.Cells.Item(566, 13).Value = 25000
$ws.Cells.Item(566, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(566, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(566, 16).Value = 1667
$ws.Cells.Item(566, 17).Value = 15
$ws.Cells.Item(566, 18).Value = "Hortaliza"

# Populate the second new row (567) - Pimiento, Zafiro verde
$ws.Cells.Item(567, 1).Value = 5
$ws.Cells.Item(567, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(567, 3).Value = "Maule"
$ws.Cells.Item(567, 4).Value = 44783
$ws.Cells.Item(567, 5).Value = 7
$ws.Cells.Item(567, 6).Value = 100112002
$ws.Cells.Item(567, 7).Value = "Pimiento"
$ws.Cells.Item(567, 8).Value = "Zafiro verde"
$ws.Cells.Item(567, 9).Value = "Primera"
$ws.Cells.Item(567, 10).Value = 300
$ws.Cells.Item(567, 11).Value = 25000
$ws.Cells.Item(567, 12).Value = 25000
$ws.Cells.Item(567, 13).Value = 25000
$ws.Cells.Item(567, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(567, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(567, 16).Value = 1667
$ws.Cells.Item(567, 17).Value = 15
$ws.Cells.Item(567, 18).Value = "Hortaliza"
